# Refined metadata to be additional tab
#
# 1) Refresh the "time_taken" timestamps on the "data" sheet (column F, rows 2-45)
# 2) Add a new "metadata" worksheet after "data" summarizing the panel query

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1) Update column F ("time_taken") timestamps on the data sheet
# ---------------------------------------------------------------------------
$timestamps = @(
    "2021-10-05 14:20:55.366447",
    "2021-10-05 14:20:55.366455",
    "2021-10-05 14:20:55.366458",
    "2021-10-05 14:20:55.366460",
    "2021-10-05 14:20:55.366463",
    "2021-10-05 14:20:55.366466",
    "2021-10-05 14:20:55.366468",
    "2021-10-05 14:20:55.366471",
    "2021-10-05 14:20:55.366474",
    "2021-10-05 14:20:55.366476",
    "2021-10-05 14:20:55.366479",
    "2021-10-05 14:20:55.366481",
    "2021-10-05 14:20:55.366484",
    "2021-10-05 14:20:55.366486",
    "2021-10-05 14:20:55.366489",
    "2021-10-05 14:20:55.366492",
    "2021-10-05 14:20:55.366494",
    "2021-10-05 14:20:55.366497",
    "2021-10-05 14:20:55.366500",
    "2021-10-05 14:20:55.366502",
    "2021-10-05 14:20:55.366505",
    "2021-10-05 14:20:55.366507",
    "2021-10-05 14:20:55.366510",
    "2021-10-05 14:20:55.366512",
    "2021-10-05 14:20:55.366515",
    "2021-10-05 14:20:55.366518",
    "2021-10-05 14:20:55.366520",
    "2021-10-05 14:20:55.366523",
    "2021-10-05 14:20:55.366525",
    "2021-10-05 14:20:55.366528",
    "2021-10-05 14:20:55.366530",
    "2021-10-05 14:20:55.366533",
    "2021-10-05 14:20:55.366536",
    "2021-10-05 14:20:55.366538",
    "2021-10-05 14:20:55.366541",
    "2021-10-05 14:20:55.366544",
    "2021-10-05 14:20:55.366546",
    "2021-10-05 14:20:55.366549",
    "2021-10-05 14:20:55.366551",
    "2021-10-05 14:20:55.366554",
    "2021-10-05 14:20:55.366557",
    "2021-10-05 14:20:55.366559",
    "2021-10-05 14:20:55.366562",
    "2021-10-05 14:20:55.366564"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}

# ---------------------------------------------------------------------------
# 2) Add the "metadata" worksheet right after "data"
# ---------------------------------------------------------------------------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"
$metaSheet.Outline.SummaryRow = 1
$metaSheet.Outline.SummaryColumn = 1

# Header row (bold, bordered, centered + top-aligned) -- mirrors the "data"
# sheet's header styling.
$headerRange = $metaSheet.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# Data row 2 -- A2 carries the same bold/bordered/centered style as the
# "data" sheet's index column.
$a2 = $metaSheet.Range("A2")
$a2.Font.Bold = $true
$a2.Borders.LineStyle = 1
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Value = 0

$metaSheet.Cells.Item(2, 2).Value = "Hypogonadotropic hypogonadism idiopathic"
$metaSheet.Cells.Item(2, 3).Value = 650

# Force text (not numeric) storage for these -- they're identifiers /
# timestamps / version strings, not numbers for Excel to reformat.
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "1.46"

$metaSheet.Cells.Item(2, 5).Value = "2021-09-28T15:01:43.606024Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:20:55.363247"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/650/?format=json"

$metaSheet.Range("A1").Select()

# Keep "data" as the active/selected sheet, matching the un-touched
# bookViews/activeTab in the source workbook.
$dataSheet.Activate()
$dataSheet.Range("A1").Select()
